$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.043777964896466
$ws.Range("D2").Value = 1.046422437487553
$ws.Range("E2").Value = 1.050702238819456
$ws.Range("F2").Value = 1.059117335349212
$ws.Range("I2").Value = 1.037084195325046
$ws.Range("J2").Value = 1.048846589159147
$ws.Range("K2").Value = 1.049187917374362
$ws.Range("L2").Value = 1.053455782578758
$ws.Range("M2").Value = 1.061847714574401
$ws.Range("N2").Value = 1.020083205544521

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.045140047200246
$ws.Range("D3").Value = 1.04769971047788
$ws.Range("E3").Value = 1.051951329093779
$ws.Range("F3").Value = 1.060545101620936
$ws.Range("I3").Value = 1.037366805308733
$ws.Range("J3").Value = 1.049853613378434
$ws.Range("K3").Value = 1.050275620169715
$ws.Range("L3").Value = 1.054516247850206
$ws.Range("M3").Value = 1.063088094530384
$ws.Range("N3").Value = 1.02042292249447

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.046020028938228
$ws.Range("D4").Value = 1.048525116865004
$ws.Range("E4").Value = 1.052758722304576
$ws.Range("F4").Value = 1.061468339309538
$ws.Range("I4").Value = 1.037547084364344
$ws.Range("J4").Value = 1.050503405756327
$ws.Range("K4").Value = 1.050977814091634
$ws.Range("L4").Value = 1.055201031819742
$ws.Range("M4").Value = 1.063889558649313
$ws.Range("N4").Value = 1.020641979186582

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.046389648833521
$ws.Range("D5").Value = 1.048871864586153
$ws.Range("E5").Value = 1.053097950490324
$ws.Range("F5").Value = 1.061856324151796
$ws.Range("I5").Value = 1.037622255082367
$ws.Range("J5").Value = 1.050776146393683
$ws.Range("K5").Value = 1.051272631715084
$ws.Range("L5").Value = 1.055488581563896
$ws.Range("M5").Value = 1.064226224369151
$ws.Range("N5").Value = 1.020733889118435

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.046451690764246
$ws.Range("D6").Value = 1.048930070330933
$ws.Range("E6").Value = 1.053154896787194
$ws.Range("F6").Value = 1.061921460191014
$ws.Range("I6").Value = 1.037634840333317
$ws.Range("J6").Value = 1.050821915511962
$ws.Range("K6").Value = 1.051322110440144
$ws.Range("L6").Value = 1.05553684298343
$ws.Range("M6").Value = 1.064282736326513
$ws.Range("N6").Value = 1.020749310589915

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.046024969085805
$ws.Range("D7").Value = 1.04852975111452
$ws.Range("E7").Value = 1.052763255868835
$ws.Range("F7").Value = 1.061473524143039
$ws.Range("I7").Value = 1.037548091228502
$ws.Range("J7").Value = 1.050507051823492
$ws.Range("K7").Value = 1.050981754966174
$ws.Range("L7").Value = 1.055204875379203
$ws.Range("M7").Value = 1.06389405824803
$ws.Range("N7").Value = 1.02064320800337

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.04423857511987
$ws.Range("D8").Value = 1.046854322894172
$ws.Range("E8").Value = 1.051124553486915
$ws.Range("F8").Value = 1.059599987425546
$ws.Range("I8").Value = 1.037180241324097
$ws.Range("J8").Value = 1.049187296618041
$ws.Range("K8").Value = 1.049555849836301
$ws.Range("L8").Value = 1.053814465508547
$ws.Range("M8").Value = 1.062267146355431
$ws.Range("N8").Value = 1.020198173126561

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.041079946459546
$ws.Range("D9").Value = 1.043893586557746
$ws.Range("E9").Value = 1.048230247647196
$ws.Range("F9").Value = 1.056293593806342
$ws.Range("I9").Value = 1.036512167744824
$ws.Range("J9").Value = 1.046847622973729
$ws.Range("K9").Value = 1.04703062764189
$ws.Range("L9").Value = 1.051353438236902
$ws.Range("M9").Value = 1.059391362008001
$ws.Range("N9").Value = 1.01940807063616

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.038966623650252
$ws.Range("D10").Value = 1.041913837508065
$ws.Range("E10").Value = 1.046295949384131
$ws.Range("F10").Value = 1.054085678443585
$ws.Range("I10").Value = 1.036053349724307
$ws.Range("J10").Value = 1.0452781357979
$ws.Range("K10").Value = 1.04533843885573
$ws.Range("L10").Value = 1.04970516061992
$ws.Range("M10").Value = 1.057467879582648
$ws.Range("N10").Value = 1.018877299793899

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.038049662407177
$ws.Range("D11").Value = 1.041055118996113
$ws.Range("E11").Value = 1.045457190082369
$ws.Range("F11").Value = 1.053128691953644
$ws.Range("I11").Value = 1.035851472727108
$ws.Range("J11").Value = 1.044596177176998
$ws.Range("K11").Value = 1.044603582995271
$ws.Range("L11").Value = 1.048989584724927
$ws.Range("M11").Value = 1.056633436836252
$ws.Range("N11").Value = 1.018646496331948

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.037708773500474
$ws.Range("D12").Value = 1.040735925951573
$ws.Range("E12").Value = 1.045145452874526
$ws.Range("F12").Value = 1.052773076660283
$ws.Range("I12").Value = 1.035776003428113
$ws.Range("J12").Value = 1.044342508216358
$ws.Range("K12").Value = 1.044330300475684
$ws.Range("L12").Value = 1.048723503774488
$ws.Range("M12").Value = 1.056323247583272
$ws.Range("N12").Value = 1.018560617516866

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.037781908495281
$ws.Range("D13").Value = 1.04080440428967
$ws.Range("E13").Value = 1.045212330001732
$ws.Range("F13").Value = 1.052849364137416
$ws.Range("I13").Value = 1.035792213739819
$ws.Range("J13").Value = 1.044396937396039
$ws.Range("K13").Value = 1.044388935257864
$ws.Range("L13").Value = 1.048780591977207
$ws.Range("M13").Value = 1.056389795181012
$ws.Range("N13").Value = 1.018579045541325

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.038021490350137
$ws.Range("D14").Value = 1.041028739055959
$ws.Range("E14").Value = 1.045431425583926
$ws.Range("F14").Value = 1.053099299734527
$ws.Range("I14").Value = 1.035845244277445
$ws.Range("J14").Value = 1.044575216174925
$ws.Range("K14").Value = 1.044581000021088
$ws.Range("L14").Value = 1.048967596215162
$ws.Range("M14").Value = 1.056607801424261
$ws.Range("N14").Value = 1.018639400592732

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.038169066192339
$ws.Range("D15").Value = 1.041166928903284
$ws.Range("E15").Value = 1.045566392982052
$ws.Range("F15").Value = 1.053253273526442
$ws.Range("I15").Value = 1.035877854066872
$ws.Range("J15").Value = 1.044685011884778
$ws.Range("K15").Value = 1.044699294339444
$ws.Range("L15").Value = 1.049082777899159
$ws.Range("M15").Value = 1.056742090293806
$ws.Range("N15").Value = 1.018676567658569

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.039027439220669
$ws.Range("D16").Value = 1.041970796319647
$ws.Range("E16").Value = 1.046351589459911
$ws.Range("F16").Value = 1.054149170073989
$ws.Range("I16").Value = 1.036066679961024
$ws.Range("J16").Value = 1.045323345028725
$ws.Range("K16").Value = 1.045387163599339
$ws.Range("L16").Value = 1.049752611435448
$ws.Range("M16").Value = 1.057523225429586
$ws.Range("N16").Value = 1.018892596767108

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.039565366474757
$ws.Range("D17").Value = 1.042474642733425
$ws.Range("E17").Value = 1.046843798864643
$ws.Range("F17").Value = 1.054710885206912
$ws.Range("I17").Value = 1.036184266084512
$ws.Range("J17").Value = 1.045723119419802
$ws.Range("K17").Value = 1.045818073004561
$ws.Range("L17").Value = 1.050172278863627
$ws.Range("M17").Value = 1.058012788852964
$ws.Range("N17").Value = 1.019027843731992

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.039878949530133
$ws.Range("D18").Value = 1.042768385952584
$ws.Range("E18").Value = 1.047130781290427
$ws.Range("F18").Value = 1.055038433411512
$ws.Range("I18").Value = 1.036252542878412
$ws.Range("J18").Value = 1.045956073632843
$ws.Range("K18").Value = 1.046069210276792
$ws.Range("L18").Value = 1.05041688454927
$ws.Range("M18").Value = 1.058298192554708
$ws.Range("N18").Value = 1.019106636911218

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.039985842773493
$ws.Range("D19").Value = 1.042868520894578
$ws.Range("E19").Value = 1.047228615528747
$ws.Range("F19").Value = 1.055150103675562
$ws.Range("I19").Value = 1.036275771101302
$ws.Range("J19").Value = 1.046035466560838
$ws.Range("K19").Value = 1.046154806986844
$ws.Range("L19").Value = 1.05050025850652
$ws.Range("M19").Value = 1.058395482499059
$ws.Range("N19").Value = 1.019133487434511

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.039507670721701
$ws.Range("D20").Value = 1.042420599495539
$ws.Range("E20").Value = 1.046791001386002
$ws.Range("F20").Value = 1.054650627918636
$ws.Range("I20").Value = 1.036171682200413
$ws.Range("J20").Value = 1.04568025096899
$ws.Range("K20").Value = 1.045771861732502
$ws.Range("L20").Value = 1.050127271093193
$ws.Range("M20").Value = 1.057960278923578
$ws.Range("N20").Value = 1.019013342755593

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.037950947486133
$ws.Range("D21").Value = 1.040962684360915
$ws.Range("E21").Value = 1.045366912534115
$ws.Range("F21").Value = 1.053025704061805
$ws.Range("I21").Value = 1.035829641457049
$ws.Range("J21").Value = 1.044522727469508
$ws.Range("K21").Value = 1.044524450727234
$ws.Range("L21").Value = 1.048912536010846
$ws.Range("M21").Value = 1.05654361068953
$ws.Range("N21").Value = 1.018621631635481

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.036970500396627
$ws.Range("D22").Value = 1.040044720666069
$ws.Range("E22").Value = 1.044470460440993
$ws.Range("F22").Value = 1.052003191227852
$ws.Range("I22").Value = 1.035611790465191
$ws.Range("J22").Value = 1.043792866201747
$ws.Range("K22").Value = 1.043738275023962
$ws.Range("L22").Value = 1.048147137221162
$ws.Range("M22").Value = 1.055651503520604
$ws.Range("N22").Value = 1.018374489517879

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.037490414734166
$ws.Range("D23").Value = 1.040531476973249
$ws.Range("E23").Value = 1.044945789868182
$ws.Range("F23").Value = 1.052545328035589
$ws.Range("I23").Value = 1.035727542992246
$ws.Range("J23").Value = 1.044179978188622
$ws.Range("K23").Value = 1.04415522134461
$ws.Range("L23").Value = 1.048553047240905
$ws.Range("M23").Value = 1.056124560122454
$ws.Range("N23").Value = 1.01850558603382

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.039533741507248
$ws.Range("D24").Value = 1.042445019746498
$ws.Range("E24").Value = 1.046814858648955
$ws.Range("F24").Value = 1.054677855873414
$ws.Range("I24").Value = 1.036177369271268
$ws.Range("J24").Value = 1.045699622079071
$ws.Range("K24").Value = 1.045792743249982
$ws.Range("L24").Value = 1.050147608720791
$ws.Range("M24").Value = 1.057984006364643
$ws.Range("N24").Value = 1.019019895413266

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.041897837920476
$ws.Range("D25").Value = 1.044660032567621
$ws.Range("E25").Value = 1.048979315460049
$ws.Range("F25").Value = 1.057148997308302
$ws.Range("I25").Value = 1.036687243481574
$ws.Range("J25").Value = 1.047454178429282
$ws.Range("K25").Value = 1.047684975167296
$ws.Range("L25").Value = 1.051990993892163
$ws.Range("M25").Value = 1.060135910404223
$ws.Range("N25").Value = 1.019613037023756
